$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0214 to SCD0012
$ws.Name = "SCD0012"

# Update TC_ID cell B2 from "DGS-229" to "SCD0012-004"
$ws.Range("B2").Value = "SCD0012-004"

# Column B widened to fit the longer TC_ID text ("SCD0012-004")
$ws.Range("B:B").ColumnWidth = 11.6

# Update selection to reflect the user's last interaction (B3)
$ws.Range("B3").Select()
